$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.813.59'
$ws.Range('E2').Value = '  +2.39%  '
$ws.Range('D3').Value = '3.048.72'
$ws.Range('E3').Value = '  +2.67%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '525.17'
$ws.Range('E5').Value = '  +5.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.34'
$ws.Range('E6').Value = '  +5.90%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +5.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.62'
$ws.Range('E9').Value = '  +4.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.112'
$ws.Range('E10').Value = '  +8.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.368'
$ws.Range('E11').Value = '  +5.36%  '
$ws.Range('E12').Value = '  +2.75%  '
$ws.Range('D13').Value = '3.574.27'
$ws.Range('E13').Value = '  +2.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.59'
$ws.Range('E14').Value = '  +7.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000169'
$ws.Range('E15').Value = '  +16.69%  '
$ws.Range('D16').Value = '57.798.08'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.19'
$ws.Range('E17').Value = '  +7.09%  '
$ws.Range('D18').Value = '3.052.29'
$ws.Range('E18').Value = '  +2.74%  '
$ws.Range('E19').Value = '  +6.04%  '
$ws.Range('E20').Value = '  +6.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '340.39'
$ws.Range('E21').Value = '  +5.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.499'
$ws.Range('E23').Value = '  +8.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.51'
$ws.Range('E24').Value = '  +7.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.173'
$ws.Range('E25').Value = '  +7.51%  '
$ws.Range('D26').Value = '0.0₃0965'
$ws.Range('E26').Value = '  +8.59%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.93'
$ws.Range('E28').Value = '  +7.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.32'
$ws.Range('E29').Value = '  +8.84%  '
$ws.Range('E30').Value = '  +7.72%  '
$ws.Range('E31').Value = '  +6.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.03'
$ws.Range('E32').Value = '  +6.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '157.98'
$ws.Range('E33').Value = '  +1.87%  '
$ws.Range('E34').Value = '  +6.51%  '
$ws.Range('E35').Value = '  +5.93%  '
$ws.Range('E36').Value = '  +4.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '25.86'
$ws.Range('E37').Value = '  +11.62%  '
$ws.Range('E38').Value = '  +4.13%  '
$ws.Range('D39').Value = '3.083.32'
$ws.Range('E39').Value = '  +2.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.71'
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.85'
$ws.Range('E41').Value = '  +8.57%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.663'
$ws.Range('E43').Value = '  +4.86%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.332.23'
$ws.Range('E44').Value = '  +6.27%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.47'
$ws.Range('E45').Value = '  +5.21%  '
$ws.Range('E46').Value = '  +3.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.05'
$ws.Range('E47').Value = '  +5.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0247'
$ws.Range('E48').Value = '  +5.23%  '
$ws.Range('E49').Value = '  +6.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.98'
$ws.Range('E50').Value = '  +4.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0894'
$ws.Range('E51').Value = '  +5.64%  '
